# Scripting SCD0331 until SCD0333
# Update the TC_ID in row 2 from DGS-347 to DGS-346 and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "DGS-346"

$ws.Range("B3").Select()
